$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.847.49"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'240.69"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'0.6282"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.07689"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "'0.2919"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'24.75"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.07736"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "1.843.19"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'5.026"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "'0.6797"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "'83.49"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "'6.170"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "29.447.34"
$ws.Range("D19").Value = "'227.82"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'12.41"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'7.404"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'157.50"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'0.1374"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "'8.400"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'17.67"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "'1.350"
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'4.115"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'1.160"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'0.7076"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "1.226.30"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D40").Value = "'6.528"
$ws.Range("E40").Value = "  +4.52%  "
$ws.Range("D41").Value = "'0.9069"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'66.09"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").Value = "'7.158"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "'0.4015"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'8.996"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").Value = "'1.672"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'0.05714"
$ws.Range("E51").Value = "  +0.11%  "
